$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample number text from E7420 to E7420L for all rows (G2:G27)
# so the shared string table entry itself is updated in place.
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("G$r").Value = "E7420L"
}

# Convert the H2:H27 "FALSE()" formula cells into literal boolean FALSE values.
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("H$r").Value = $false
}
